$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data rows (2-6) so the shared-string table rebuilds in column-fill order
# matching how the full 23-row dataset (rows 2-24) was originally authored.
$ws.Range("A2:F6").Clear()

# Column B (daytime)
$ws.Range("B2").Value = "11-11-21"
$ws.Range("B3").Value = "13-11-21"
$ws.Range("B4").Value = "13-11-21"
$ws.Range("B5").Value = "13-11-21"
$ws.Range("B6").Value = "13-11-21"
$ws.Range("B7").Value = "13-11-21"
$ws.Range("B8").Value = "13-11-21"
$ws.Range("B9").Value = "13-11-21"
$ws.Range("B10").Value = "15-11-21"
$ws.Range("B11").Value = "15-11-21"
$ws.Range("B12").Value = "15-11-21"
$ws.Range("B13").Value = "15-11-21"
$ws.Range("B14").Value = "15-11-21"
$ws.Range("B15").Value = "15-11-21"
$ws.Range("B16").Value = "15-11-21"
$ws.Range("B17").Value = "15-11-21"
$ws.Range("B18").Value = "16-11-21"
$ws.Range("B19").Value = "16-11-21"
$ws.Range("B20").Value = "16-11-21"
$ws.Range("B21").Value = "16-11-21"
$ws.Range("B22").Value = "16-11-21"
$ws.Range("B23").Value = "16-11-21"
$ws.Range("B24").Value = "16-11-21"

# Column C (time)
$ws.Range("C2").Value = "14:11"
$ws.Range("C3").Value = "19:15"
$ws.Range("C4").Value = "19:15"
$ws.Range("C5").Value = "19:15"
$ws.Range("C6").Value = "19:16"
$ws.Range("C7").Value = "23:42"
$ws.Range("C8").Value = "23:42"
$ws.Range("C9").Value = "23:42"
$ws.Range("C10").Value = "16:57"
$ws.Range("C11").Value = "17:04"
$ws.Range("C12").Value = "17:13"
$ws.Range("C13").Value = "17:13"
$ws.Range("C14").Value = "17:13"
$ws.Range("C15").Value = "17:15"
$ws.Range("C16").Value = "17:58"
$ws.Range("C17").Value = "19:10"
$ws.Range("C18").Value = "11:16"
$ws.Range("C19").Value = "11:16"
$ws.Range("C20").Value = "11:23"
$ws.Range("C21").Value = "11:27"
$ws.Range("C22").Value = "11:28"
$ws.Range("C23").Value = "11:29"
$ws.Range("C24").Value = "11:35"

# Column D (codfisc)
$ws.Range("D2").Value = "jhbhb"
$ws.Range("D3").Value = "uihiuh"
$ws.Range("D4").Value = "uyuy"
$ws.Range("D5").Value = "ctct"
$ws.Range("D6").Value = "yvy"
$ws.Range("D7").Value = "dsiofsdif"
$ws.Range("D8").Value = "uggu"
$ws.Range("D9").Value = "derf"
$ws.Range("D10").Value = "edewd"
$ws.Range("D11").Value = "sdasssdsadsdsadasdsad"
$ws.Range("D12").Value = "sdd"
$ws.Range("D13").Value = "sdsad"
$ws.Range("D14").Value = "dasdsad"
$ws.Range("D15").Value = "xczxc"
$ws.Range("D16").Value = "sadasd"
$ws.Range("D17").Value = "dsdsad"
$ws.Range("D18").Value = "CNTFBA75P24H501P"
$ws.Range("D19").Value2 = "1.603008416904999E+26"
$ws.Range("D20").Value = "CNTFBA75P24H501P"
$ws.Range("D21").Value = "CNTFBA75P24H501P"
$ws.Range("D22").Value = "CNTFBA75P24H501P"
$ws.Range("D23").Value = "CNTFBA75P24H501P"
$ws.Range("D24").Value = "dss"

# Column E (ticket)
$ws.Range("E2").Value = "dxd"
$ws.Range("E3").Value = "ihhi"
$ws.Range("E4").Value = "ttrct"
$ws.Range("E5").Value = "crtct"
$ws.Range("E6").Value = "vyyvy"
$ws.Range("E7").Value = "iniunin"
$ws.Range("E8").Value = "guybub"
$ws.Range("E9").Value = "crct"
$ws.Range("E10").Value = "ewede"
$ws.Range("E11").Value = "sadsadsadsadsadasdsadsad"
$ws.Range("E12").Value = "dsads"
$ws.Range("E13").Value = "sdsads"
$ws.Range("E14").Value = "sadsadsad"
$ws.Range("E15").Value = "czxcxzc"
$ws.Range("E16").Value = "sdasd"
$ws.Range("E17").Value = "sadsadas"
$ws.Range("E18").Value2 = "99914130540219"
$ws.Range("E19").Value2 = "99914130540219"
$ws.Range("E20").Value2 = "99914130540219"
$ws.Range("E21").Value2 = "1603008416904"
$ws.Range("E22").Value2 = "1603008416904"
$ws.Range("E23").Value2 = "99914130540219"
$ws.Range("E24").Value = "sdasd"

# Column A (id) and F (numIngressi)
$ws.Range("A2").Value = 4
$ws.Range("F2").Value = 5
$ws.Range("A3").Value = 6
$ws.Range("F3").Value = 0
$ws.Range("A4").Value = 7
$ws.Range("F4").Value = 5
$ws.Range("A5").Value = 8
$ws.Range("F5").Value = 1
$ws.Range("A6").Value = 9
$ws.Range("F6").Value = 4
$ws.Range("A7").Value = 10
$ws.Range("F7").Value = 8
$ws.Range("A8").Value = 11
$ws.Range("F8").Value = 9
$ws.Range("A9").Value = 12
$ws.Range("F9").Value = 6
$ws.Range("A10").Value = 13
$ws.Range("F10").Value = 1
$ws.Range("A11").Value = 14
$ws.Range("F11").Value = 1
$ws.Range("A12").Value = 15
$ws.Range("F12").Value = 1
$ws.Range("A13").Value = 16
$ws.Range("F13").Value = 1
$ws.Range("A14").Value = 17
$ws.Range("F14").Value = 1
$ws.Range("A15").Value = 18
$ws.Range("F15").Value = 1
$ws.Range("A16").Value = 19
$ws.Range("F16").Value = 1
$ws.Range("A17").Value = 20
$ws.Range("F17").Value = 1
$ws.Range("A18").Value = 24
$ws.Range("F18").Value = 8
$ws.Range("A19").Value = 25
$ws.Range("F19").Value = 1
$ws.Range("A20").Value = 26
$ws.Range("F20").Value = 5
$ws.Range("A21").Value = 27
$ws.Range("F21").Value = 5
$ws.Range("A22").Value = 29
$ws.Range("F22").Value = 5
$ws.Range("A23").Value = 30
$ws.Range("F23").Value = 3
$ws.Range("A24").Value = 32
$ws.Range("F24").Value = 1

Write-Host "done"
